$d = $word.ActiveDocument

$d.Content.Find.Execute("85×13=", $false, $false, $false, $false, $false, $true, 1, $false, "48×74=", 2) | Out-Null
$d.Content.Find.Execute("86×64=", $false, $false, $false, $false, $false, $true, 1, $false, "13×85=", 2) | Out-Null
$d.Content.Find.Execute("79×80=", $false, $false, $false, $false, $false, $true, 1, $false, "47×99=", 2) | Out-Null
$d.Content.Find.Execute("33×73=", $false, $false, $false, $false, $false, $true, 1, $false, "35×66=", 2) | Out-Null
$d.Content.Find.Execute("54×91=", $false, $false, $false, $false, $false, $true, 1, $false, "23×65=", 2) | Out-Null
$d.Content.Find.Execute("65×90=", $false, $false, $false, $false, $false, $true, 1, $false, "14×32=", 2) | Out-Null
$d.Content.Find.Execute("28×55=", $false, $false, $false, $false, $false, $true, 1, $false, "78×38=", 2) | Out-Null
$d.Content.Find.Execute("63×71=", $false, $false, $false, $false, $false, $true, 1, $false, "21×50=", 2) | Out-Null
$d.Content.Find.Execute("23×90=", $false, $false, $false, $false, $false, $true, 1, $false, "47×84=", 2) | Out-Null
$d.Content.Find.Execute("79×33=", $false, $false, $false, $false, $false, $true, 1, $false, "28×17=", 2) | Out-Null
$d.Content.Find.Execute("71×43=", $false, $false, $false, $false, $false, $true, 1, $false, "45×63=", 2) | Out-Null
$d.Content.Find.Execute("80×86=", $false, $false, $false, $false, $false, $true, 1, $false, "69×49=", 2) | Out-Null
$d.Content.Find.Execute("82×93=", $false, $false, $false, $false, $false, $true, 1, $false, "62×98=", 2) | Out-Null
$d.Content.Find.Execute("23×25=", $false, $false, $false, $false, $false, $true, 1, $false, "55×77=", 2) | Out-Null
$d.Content.Find.Execute("49×58=", $false, $false, $false, $false, $false, $true, 1, $false, "28×59=", 2) | Out-Null
$d.Content.Find.Execute("76×27=", $false, $false, $false, $false, $false, $true, 1, $false, "50×36=", 2) | Out-Null
$d.Content.Find.Execute("67×64=", $false, $false, $false, $false, $false, $true, 1, $false, "27×52=", 2) | Out-Null
$d.Content.Find.Execute("16×35=", $false, $false, $false, $false, $false, $true, 1, $false, "16×67=", 2) | Out-Null
$d.Content.Find.Execute("87×95=", $false, $false, $false, $false, $false, $true, 1, $false, "43×75=", 2) | Out-Null
$d.Content.Find.Execute("84×14=", $false, $false, $false, $false, $false, $true, 1, $false, "69×35=", 2) | Out-Null
$d.Content.Find.Execute("61×18=", $false, $false, $false, $false, $false, $true, 1, $false, "23×19=", 2) | Out-Null
$d.Content.Find.Execute("48×61=", $false, $false, $false, $false, $false, $true, 1, $false, "12×83=", 2) | Out-Null
$d.Content.Find.Execute("37×58=", $false, $false, $false, $false, $false, $true, 1, $false, "43×75=", 2) | Out-Null
$d.Content.Find.Execute("56×25=", $false, $false, $false, $false, $false, $true, 1, $false, "84×67=", 2) | Out-Null
$d.Content.Find.Execute("80×89=", $false, $false, $false, $false, $false, $true, 1, $false, "28×29=", 2) | Out-Null
